$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$officesText = "4.2% CR/LFM+CDL/H:2/Offices" + $nl +
  "24.0% CR/LFM+CDL/HBET:3-5/Offices" + $nl +
  "6.5% MR/LWAL+CDL/H:1/Offices" + $nl +
  "12.9% MR/LWAL+CDL/H:2/Offices" + $nl +
  "45.2% MR/LWAL+CDL/HBET:3-5/Offices" + $nl +
  "5.2% S/LFM+CDL/HBET:3-5/Offices" + $nl +
  "2.0% S/LWAL+CDL/HBET:3-5/Offices" + $nl +
  "0.0% CR/LFM+CDL/H:1/Offices" + $nl +
  "0.0% S/LFM+CDL/H:1/Offices" + $nl +
  "0.0% S/LFM+CDL/H:2/Offices" + $nl +
  "0.0% S/LWAL+CDL/H:1/Offices" + $nl +
  "0.0% S/LWAL+CDL/H:2/Offices" + $nl +
  "0.0% W/LFM+CDL/H:1/Offices" + $nl +
  "0.0% W/LFM+CDL/H:2/Offices"

$tradeText = "2.4% CR/LFM+CDL/H:2/Trade" + $nl +
  "0.0% CR/LFM+CDL/HBET:3-5/Trade" + $nl +
  "58.2% MR/LWAL+CDL/H:1/Trade" + $nl +
  "6.5% MR/LWAL+CDL/H:2/Trade" + $nl +
  "0.0% MR/LWAL+CDL/HBET:3-5/Trade" + $nl +
  "0.0% S/LFM+CDL/HBET:3-5/Trade" + $nl +
  "0.0% S/LWAL+CDL/HBET:3-5/Trade" + $nl +
  "21.6% CR/LFM+CDL/H:1/Trade" + $nl +
  "4.6% S/LFM+CDL/H:1/Trade" + $nl +
  "0.5% S/LFM+CDL/H:2/Trade" + $nl +
  "1.8% S/LWAL+CDL/H:1/Trade" + $nl +
  "0.2% S/LWAL+CDL/H:2/Trade" + $nl +
  "3.8% W/LFM+CDL/H:1/Trade" + $nl +
  "0.4% W/LFM+CDL/H:2/Trade"

$hotelsText = "2.4% CR/LFM+CDL/H:2/Hotels" + $nl +
  " 16.8% CR/LFM+CDL/HBET:3-5/Hotels" + $nl +
  " 12.9% MR/LWAL+CDL/H:1/Hotels" + $nl +
  " 6.5% MR/LWAL+CDL/H:2/Hotels" + $nl +
  " 45.3% MR/LWAL+CDL/HBET:3-5/Hotels" + $nl +
  " 5.2% S/LFM+CDL/HBET:3-5/Hotels" + $nl +
  " 2.0% S/LWAL+CDL/HBET:3-5/Hotels" + $nl +
  " 4.8% CR/LFM+CDL/H:1/Hotels" + $nl +
  " 0.0% S/LFM+CDL/H:1/Hotels" + $nl +
  " 0.0% S/LFM+CDL/H:2/Hotels" + $nl +
  " 0.0% S/LWAL+CDL/H:1/Hotels" + $nl +
  " 0.0% S/LWAL+CDL/H:2/Hotels" + $nl +
  " 0.8% W/LFM+CDL/H:1/Hotels" + $nl +
  " 3.3% W/LFM+CDL/H:2/Hotels"

$ws.Range("B2:D2").Style = "Normal"

$ws.Range("B2").Value = $officesText
$ws.Range("C2").Value = $tradeText
$ws.Range("D2").Value = $hotelsText

$ws.Columns.Item(2).ColumnWidth = 8.83203125
$ws.Columns.Item(3).ColumnWidth = 8.83203125
$ws.Columns.Item(4).ColumnWidth = 8.83203125

$ws.Rows.Item(2).EntireRow.AutoFit()
